$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.02,1.02180215628182,1.027618200386493,0.9926147277508489,1.020239311050894,1.030946282868828,1.026991572615026,1.03043740429647,0.9955398523336033,1.023080185139902,1.013073959379771)
    3 = @(1.02,1.022577510827691,1.028189888896194,0.9936372048519304,1.021650225594853,1.031088335780921,1.027405608717223,1.030817494064643,0.9963617723202692,1.024295637257594,1.013211276150165)
    4 = @(1.02,1.02307936737259,1.028559895784052,0.9942998659930995,1.022563404435056,1.031179123125931,1.027673003403625,1.031062822506086,0.9968940712668345,1.025081830657428,1.013299941713283)
    5 = @(1.02,1.02329038286283,1.028715466100335,0.9945786998346017,1.022947361226487,1.031217019126009,1.027785292441931,1.03116581063787,0.997117960005301,1.025412280566139,1.013337171599949)
    6 = @(1.02,1.023325815315823,1.028741588144393,0.9946255319796338,1.023011832658306,1.031223366128491,1.027804138996416,1.031183094102401,0.9971555583673453,1.025467760775004,1.013343420009939)
    7 = @(1.02,1.023082186835614,1.028561974447459,0.9943035907982488,1.022568534659017,1.03117963055946,1.027674504301953,1.031064199220749,0.9968970624462087,1.025086246400664,1.013300439358393)
    8 = @(1.02,1.022064159267292,1.027811386619518,0.9929600610674301,1.020716091946442,1.030994524139968,1.027131604225888,1.030565984344612,0.995817528259106,1.023491013488055,1.013120404935929)
    9 = @(1.02,1.020271458995708,1.02648946843368,0.9906006454969559,1.0174534199693,1.030659702358139,1.026171029548005,1.029683386142633,0.9939188001724441,1.020677710263055,1.012801733962393)
    10 = @(1.02,1.019077185305457,1.025608740399385,0.989033133672735,1.015279166160612,1.030430700518498,1.025528045166171,1.02909188123661,0.9926553831429383,1.018800468999474,1.012588339087986)
    11 = @(1.02,1.018560267062506,1.025227520798733,0.988355674866747,1.014337850764017,1.030330172011376,1.025249015062506,1.028835024584791,0.9921088820399291,1.017987162046739,1.012495714433539)
    12 = @(1.02,1.018368293057168,1.025085941277826,0.9881042295826724,1.013988223137411,1.030292625869322,1.025145279121901,1.028739507368028,0.9919059725120875,1.017684992961183,1.012461276150124)
    13 = @(1.02,1.018409470650747,1.025116309532878,0.9881581567098651,1.014063218660085,1.030300688932792,1.025167534977315,1.02876000106772,0.9919494934313052,1.017749812505212,1.012468664785668)
    14 = @(1.02,1.018544397764411,1.02521581732747,0.9883348863814464,1.014308950073575,1.030327072623324,1.025240442086077,1.028827131326315,0.9920921077337197,1.017962186116056,1.012492868433839)
    15 = @(1.02,1.018627535137771,1.025277130368384,0.9884438009545853,1.014460355686917,1.030343301281,1.025285350419584,1.028868478040118,0.9921799884222134,1.018093027046533,1.012507776684798)
    16 = @(1.02,1.01911149596901,1.025634043768248,0.9890781214508737,1.015341641014013,1.030437343445007,1.025546550585657,1.02910891261601,0.9926916645766087,1.018854435768494,1.012594481592741)
    17 = @(1.02,1.019415128692437,1.025857964806815,0.989476357848556,1.015894485696215,1.030495967268323,1.025710230586125,1.029259535571963,0.9930127773699352,1.019331924852089,1.012648809644151)
    18 = @(1.02,1.019592252740597,1.025988587730579,0.9897087662937556,1.016216965458812,1.030530029422765,1.025805643125185,1.029347320809367,0.9932001317071769,1.019610393039569,1.012680476723466)
    19 = @(1.02,1.019652650866345,1.02603312903832,0.9897880325774034,1.016326925353438,1.030541621315671,1.025838166274018,1.029377241315371,0.9932640239640975,1.019705336279145,1.012691270717729)
    20 = @(1.02,1.01938254967719,1.025833938786466,0.9894336180360679,1.015835169189023,1.030489691151312,1.025692675394853,1.02924338245061,0.9929783193494215,1.019280699274802,1.012642982986492)
    21 = @(1.02,1.018504664189741,1.025186514134142,0.9882828385668249,1.014236587808246,1.030319308950876,1.025218975288453,1.028807366148349,0.9920501090198102,1.01789964935752,1.012485741982243)
    22 = @(1.02,1.017952890892586,1.024779582800851,0.9875604150241495,1.013231603247971,1.030210994804541,1.02492061069455,1.028532593787705,0.9914670000341481,1.017030917992433,1.012386685384294)
    23 = @(1.02,1.018245378231797,1.024995292068762,0.9879432794643023,1.013764355720573,1.030268526691074,1.025078829549183,1.028678315476793,0.991776070289318,1.017491488904011,1.012439215389479)
    24 = @(1.02,1.019397270669409,1.025844795069053,0.9894529299347244,1.015861971718667,1.03049252746671,1.025700608013166,1.029250681568975,0.9929938892766442,1.019303846040616,1.012645615868885)
    25 = @(1.02,1.020734767582624,1.026831124073653,0.9912096547607049,1.018296732923734,1.030747283488224,1.026419821889978,1.029912109351004,0.9944092447426414,1.021405305491688,1.012884285940956)
}

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item([int]$row, $cols[$i]).Value = $vals[$i]
    }
}

Write-Output "applied"